# Generate Report for Archive
$wb = $excel.ActiveWorkbook

# 1. Update status text: "Ready for handoff" -> "In Translation"
#    Overview sheet: the zh-cn (E2) and de-de (F2) status cells
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

#    zh-cn / de-de sheets: the Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# 2. Narrow the language status columns (zh-cn / de-de) on the Overview sheet
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# 3. Narrow the "Status" column on the per-locale sheets
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
